$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "314.48"
Set-TextValue "E2" "3.18%"
Set-TextValue "G2" "14"
Set-TextValue "D3" "34.89"
Set-TextValue "E3" "-2.26%"
Set-TextValue "G3" "14"
Set-TextValue "D4" "5.117"
Set-TextValue "E4" "1.04%"
Set-TextValue "G4" "14"
Set-TextValue "D5" "0.08155"
Set-TextValue "E5" "3.38%"
Set-TextValue "G5" "14"
Set-TextValue "D6" "2.131"
Set-TextValue "E6" "1.21%"
Set-TextValue "G6" "14"
Set-TextValue "D7" "4.144"
Set-TextValue "E7" "0.34%"
Set-TextValue "G7" "14"
Set-TextValue "D8" "7.957"
Set-TextValue "E8" "0.47%"
Set-TextValue "G8" "14"
Set-TextValue "D9" "0.9336"
Set-TextValue "E9" "1.28%"
Set-TextValue "G9" "14"
Set-TextValue "D10" "0.1043"
Set-TextValue "E10" "6.93%"
Set-TextValue "G10" "14"
Set-TextValue "E11" "5.08%"
Set-TextValue "G11" "14"
Set-TextValue "D12" "0.09090"
Set-TextValue "E12" "4.64%"
Set-TextValue "G12" "14"
Set-TextValue "D13" "0.03635"
Set-TextValue "E13" "2.66%"
Set-TextValue "G13" "14"
Set-TextValue "D14" "0.09886"
Set-TextValue "E14" "-0.45%"
Set-TextValue "G14" "14"
Set-TextValue "D15" "0.001432"
Set-TextValue "E15" "-2.35%"
Set-TextValue "G15" "14"
Set-TextValue "D16" "0.005774"
Set-TextValue "E16" "0.94%"
Set-TextValue "G16" "14"
Set-TextValue "D17" "3.469"
Set-TextValue "E17" "0.38%"
Set-TextValue "G17" "14"
Set-TextValue "D18" "2.813"
Set-TextValue "E18" "2.29%"
Set-TextValue "G18" "14"
Set-TextValue "E19" "1.08%"
Set-TextValue "G19" "14"
Set-TextValue "D20" "0.1332"
Set-TextValue "E20" "-0.72%"
Set-TextValue "G20" "14"
Set-TextValue "D21" "5.094"
Set-TextValue "E21" "-1.57%"
Set-TextValue "G21" "14"
Set-TextValue "D22" "0.2215"
Set-TextValue "E22" "0.02%"
Set-TextValue "G22" "14"
Set-TextValue "D23" "0.04565"
Set-TextValue "E23" "1.47%"
Set-TextValue "G23" "14"
Set-TextValue "D24" "0.001248"
Set-TextValue "E24" "0.61%"
Set-TextValue "G24" "14"
Set-TextValue "D25" "0.004692"
Set-TextValue "E25" "-3.43%"
Set-TextValue "G25" "14"
Set-TextValue "D26" "0.0001251"
Set-TextValue "E26" "-3.97%"
Set-TextValue "G26" "14"
Set-TextValue "E27" "-5.42%"
Set-TextValue "G27" "14"
Set-TextValue "G28" "14"
Set-TextValue "G29" "14"
Set-TextValue "G30" "14"
Set-TextValue "G31" "14"
Set-TextValue "G32" "14"
Set-TextValue "G33" "14"
Set-TextValue "G34" "14"
Set-TextValue "G35" "14"
Set-TextValue "G36" "14"
Set-TextValue "G37" "14"
Set-TextValue "G38" "14"
Set-TextValue "D39" "0.01956"
Set-TextValue "E39" "6.28%"
Set-TextValue "G39" "14"
Set-TextValue "D40" "0.04896"
Set-TextValue "E40" "3.71%"
Set-TextValue "G40" "14"
Set-TextValue "D41" "0.007641"
Set-TextValue "E41" "-3.33%"
Set-TextValue "G41" "14"
Set-TextValue "D42" "0.1387"
Set-TextValue "E42" "-0.39%"
Set-TextValue "G42" "14"
Set-TextValue "D43" "0.007879"
Set-TextValue "E43" "1.63%"
Set-TextValue "G43" "14"
Set-TextValue "D44" "0.002142"
Set-TextValue "E44" "1.15%"
Set-TextValue "G44" "14"
Set-TextValue "D45" "0.01177"
Set-TextValue "E45" "4.57%"
Set-TextValue "G45" "14"
Set-TextValue "D46" "0.00006755"
Set-TextValue "E46" "7.21%"
Set-TextValue "G46" "14"
Set-TextValue "E47" "-0.17%"
Set-TextValue "G47" "14"
Set-TextValue "D48" "154.86"
Set-TextValue "E48" "206.21%"
Set-TextValue "G48" "14"
Set-TextValue "E49" "-10.69%"
Set-TextValue "G49" "14"
Set-TextValue "E50" "-0.17%"
Set-TextValue "G50" "14"
Set-TextValue "E51" "-0.17%"
Set-TextValue "G51" "14"
